$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: A1 = "relation", B1 = "count"
$ws.Range("A1").Value = "relation"
$ws.Range("B1").Value = "count"

# Widen column A (closest achievable value to the target stored width of 64.33203125)
$ws.Columns.Item(1).ColumnWidth = 63.5

# Move selection to F15 (as seen in the saved file after edits)
$ws.Range("F15").Select()
